# Bugfixed the naive forecaster component module
#
# The YoY component forecast "staircase" table had an off-by-one bug:
# every vintage row started writing its forecast one column too early
# (and the header/date row began one period too early), leaving a stray
# leading column of data and an extra stray forecast-origin row/column
# at the end of the sheet. This script re-lays the grid out with the
# correct date alignment and trims the now-unused trailing row/column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing forecast-origin rows (23, 24) that fall outside
# the corrected table.
$ws.Range("A23:A24").EntireRow.Delete()

# Drop the trailing date column (BA) that falls outside the corrected
# table.
$ws.Range("BA1:BA22").EntireColumn.Delete()

# Stray cells left behind by the old (buggy) column offset -- they sit
# one column to the left of where each row's corrected data now starts,
# so they must be cleared rather than overwritten.
$cellsToClear = @"
5,3
6,5
7,7
8,9
9,11
10,13
11,15
12,17
13,18
13,19
14,20
14,21
14,22
15,22
15,23
15,24
15,25
15,26
16,25
16,26
16,27
16,28
16,29
16,30
17,29
17,30
17,31
17,32
17,33
18,33
18,34
18,35
18,36
18,37
19,37
19,38
19,39
19,40
19,41
20,41
20,42
20,43
20,44
20,45
21,45
21,46
21,47
21,48
21,49
22,49
22,50
22,51
22,52
"@

foreach ($line in ($cellsToClear -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $ws.Cells.Item($r, $c).ClearContents()
}

# Corrected values (dates in row 1 / column A, and recomputed YoY
# forecast figures) for every cell of the fixed-up grid, as
# (row, column, value) triples.
$cellValues = @"
1,2,39583
1,3,39765
1,4,39948
1,5,40130
1,6,40310
1,7,40494
1,8,40676
1,9,40862
1,10,41044
1,11,41228
1,12,41409
1,13,41592
1,14,41774
1,15,41957
1,16,42137
1,17,42321
1,18,42503
1,19,42689
1,20,42867
1,21,43053
1,22,43145
1,23,43235
1,24,43326
1,25,43418
1,26,43510
1,27,43600
1,28,43691
1,29,43783
1,30,43875
1,31,43966
1,32,44068
1,33,44159
1,34,44251
1,35,44341
1,36,44432
1,37,44525
1,38,44617
1,39,44706
1,40,44798
1,41,44890
1,42,44981
1,43,45071
1,44,45163
1,45,45254
1,46,45345
1,47,45436
1,48,45534
1,49,45618
1,50,45713
1,51,45800
1,52,45891
2,1,39813
3,1,40178
3,2,-0.3422723562191532
3,3,0.4944284391569687
3,4,-0.5037688924316441
3,5,-0.5555135891318952
3,6,-0.5555135891318952
3,7,-0.5555135891318952
3,8,-0.5555135891318952
3,9,-0.5555135891318952
3,10,-0.5555135891318952
3,11,-0.5555135891318952
3,12,-0.5555135891318952
3,13,-0.5555135891318952
3,14,-0.5555135891318952
3,15,-0.5555135891318952
3,16,-0.5555135891318952
3,17,-0.5555135891318952
3,18,-0.5555135891318952
3,19,-0.5555135891318952
3,20,-0.5555135891318952
3,21,-0.5555135891318952
3,22,-0.5555135891318952
3,23,-0.5555135891318952
3,24,-0.5555135891318952
3,25,-0.5555135891318952
3,26,-0.5555135891318952
3,27,-0.5555135891318952
3,28,-0.5555135891318952
3,29,-0.5555135891318952
3,30,-0.5555135891318952
3,31,-0.5555135891318952
3,32,-0.5555135891318952
3,33,-0.5555135891318952
3,34,-0.5555135891318952
3,35,-0.5555135891318952
3,36,-0.5555135891318952
3,37,-0.5555135891318952
3,38,-0.5555135891318952
3,39,-0.5555135891318952
3,40,-0.5555135891318952
3,41,-0.5555135891318952
3,42,-0.5555135891318952
3,43,-0.5555135891318952
3,44,-0.5555135891318952
3,45,-0.5555135891318952
3,46,-0.5555135891318952
3,47,-0.5555135891318952
3,48,-0.5555135891318952
3,49,-0.5555135891318952
3,50,-0.5555135891318952
3,51,-0.5555135891318952
3,52,-0.5555135891318952
4,1,40543
4,2,-0.3678291324228367
4,3,0.3579071119161004
4,4,-0.467076459743887
4,5,-0.4782015746048418
4,6,0.169534172659791
4,7,0.8442071301477228
4,8,0.8442071301477228
4,9,0.8442071301477228
4,10,0.8442071301477228
4,11,0.8442071301477228
4,12,0.8442071301477228
4,13,0.8442071301477228
4,14,0.8442071301477228
4,15,0.8442071301477228
4,16,0.8442071301477228
4,17,0.8442071301477228
4,18,0.8442071301477228
4,19,0.8442071301477228
4,20,0.8442071301477228
4,21,0.8442071301477228
4,22,0.8442071301477228
4,23,0.8442071301477228
4,24,0.8442071301477228
4,25,0.8442071301477228
4,26,0.8442071301477228
4,27,0.8442071301477228
4,28,0.8442071301477228
4,29,0.8442071301477228
4,30,0.8442071301477228
4,31,0.8442071301477228
4,32,0.8442071301477228
4,33,0.8442071301477228
4,34,0.8442071301477228
4,35,0.8442071301477228
4,36,0.8442071301477228
4,37,0.8442071301477228
4,38,0.8442071301477228
4,39,0.8442071301477228
4,40,0.8442071301477228
4,41,0.8442071301477228
4,42,0.8442071301477228
4,43,0.8442071301477228
4,44,0.8442071301477228
4,45,0.8442071301477228
4,46,0.8442071301477228
4,47,0.8442071301477228
4,48,0.8442071301477228
4,49,0.8442071301477228
4,50,0.8442071301477228
4,51,0.8442071301477228
4,52,0.8442071301477228
5,1,40908
5,4,-0.3360947699977967
5,5,-0.3600168145845517
5,6,-0.03968684591929561
5,7,1.324233212457782
5,8,1.032338390744236
5,9,1.122475521884692
5,10,1.122475521884692
5,11,1.122475521884692
5,12,1.122475521884692
5,13,1.122475521884692
5,14,1.122475521884692
5,15,1.122475521884692
5,16,1.122475521884692
5,17,1.122475521884692
5,18,1.122475521884692
5,19,1.122475521884692
5,20,1.122475521884692
5,21,1.122475521884692
5,22,1.122475521884692
5,23,1.122475521884692
5,24,1.122475521884692
5,25,1.122475521884692
5,26,1.122475521884692
5,27,1.122475521884692
5,28,1.122475521884692
5,29,1.122475521884692
5,30,1.122475521884692
5,31,1.122475521884692
5,32,1.122475521884692
5,33,1.122475521884692
5,34,1.122475521884692
5,35,1.122475521884692
5,36,1.122475521884692
5,37,1.122475521884692
5,38,1.122475521884692
5,39,1.122475521884692
5,40,1.122475521884692
5,41,1.122475521884692
5,42,1.122475521884692
5,43,1.122475521884692
5,44,1.122475521884692
5,45,1.122475521884692
5,46,1.122475521884692
5,47,1.122475521884692
5,48,1.122475521884692
5,49,1.122475521884692
5,50,1.122475521884692
5,51,1.122475521884692
5,52,1.122475521884692
6,1,41274
6,6,-0.09284881712173787
6,7,0.6032582924254459
6,8,0.3765075513336269
6,9,0.7478380109886329
6,10,-0.1937612543835177
6,11,-0.578174579726376
6,12,-0.578174579726376
6,13,-0.578174579726376
6,14,-0.578174579726376
6,15,-0.578174579726376
6,16,-0.578174579726376
6,17,-0.578174579726376
6,18,-0.578174579726376
6,19,-0.578174579726376
6,20,-0.578174579726376
6,21,-0.578174579726376
6,22,-0.578174579726376
6,23,-0.578174579726376
6,24,-0.578174579726376
6,25,-0.578174579726376
6,26,-0.578174579726376
6,27,-0.578174579726376
6,28,-0.578174579726376
6,29,-0.578174579726376
6,30,-0.578174579726376
6,31,-0.578174579726376
6,32,-0.578174579726376
6,33,-0.578174579726376
6,34,-0.578174579726376
6,35,-0.578174579726376
6,36,-0.578174579726376
6,37,-0.578174579726376
6,38,-0.578174579726376
6,39,-0.578174579726376
6,40,-0.578174579726376
6,41,-0.578174579726376
6,42,-0.578174579726376
6,43,-0.578174579726376
6,44,-0.578174579726376
6,45,-0.578174579726376
6,46,-0.578174579726376
6,47,-0.578174579726376
6,48,-0.578174579726376
6,49,-0.578174579726376
6,50,-0.578174579726376
6,51,-0.578174579726376
6,52,-0.578174579726376
7,1,41639
7,8,0.3941427794649011
7,9,0.8822337163017169
7,10,0.6066710853121382
7,11,-0.2445716668737163
7,12,-0.8331679621937482
7,13,-0.7492845378401558
7,14,-0.7492845378401558
7,15,-0.7492845378401558
7,16,-0.7492845378401558
7,17,-0.7492845378401558
7,18,-0.7492845378401558
7,19,-0.7492845378401558
7,20,-0.7492845378401558
7,21,-0.7492845378401558
7,22,-0.7492845378401558
7,23,-0.7492845378401558
7,24,-0.7492845378401558
7,25,-0.7492845378401558
7,26,-0.7492845378401558
7,27,-0.7492845378401558
7,28,-0.7492845378401558
7,29,-0.7492845378401558
7,30,-0.7492845378401558
7,31,-0.7492845378401558
7,32,-0.7492845378401558
7,33,-0.7492845378401558
7,34,-0.7492845378401558
7,35,-0.7492845378401558
7,36,-0.7492845378401558
7,37,-0.7492845378401558
7,38,-0.7492845378401558
7,39,-0.7492845378401558
7,40,-0.7492845378401558
7,41,-0.7492845378401558
7,42,-0.7492845378401558
7,43,-0.7492845378401558
7,44,-0.7492845378401558
7,45,-0.7492845378401558
7,46,-0.7492845378401558
7,47,-0.7492845378401558
7,48,-0.7492845378401558
7,49,-0.7492845378401558
7,50,-0.7492845378401558
7,51,-0.7492845378401558
7,52,-0.7492845378401558
8,1,42004
8,10,0.7277994253700903
8,11,0.2906157075715843
8,12,-0.3970496740026364
8,13,-0.2617076051026235
8,14,0.1273541662098365
8,15,0.2751437421933511
8,16,0.2751437421933511
8,17,0.2751437421933511
8,18,0.2751437421933511
8,19,0.2751437421933511
8,20,0.2751437421933511
8,21,0.2751437421933511
8,22,0.2751437421933511
8,23,0.2751437421933511
8,24,0.2751437421933511
8,25,0.2751437421933511
8,26,0.2751437421933511
8,27,0.2751437421933511
8,28,0.2751437421933511
8,29,0.2751437421933511
8,30,0.2751437421933511
8,31,0.2751437421933511
8,32,0.2751437421933511
8,33,0.2751437421933511
8,34,0.2751437421933511
8,35,0.2751437421933511
8,36,0.2751437421933511
8,37,0.2751437421933511
8,38,0.2751437421933511
8,39,0.2751437421933511
8,40,0.2751437421933511
8,41,0.2751437421933511
8,42,0.2751437421933511
8,43,0.2751437421933511
8,44,0.2751437421933511
8,45,0.2751437421933511
8,46,0.2751437421933511
8,47,0.2751437421933511
8,48,0.2751437421933511
8,49,0.2751437421933511
8,50,0.2751437421933511
8,51,0.2751437421933511
8,52,0.2751437421933511
9,1,42369
9,12,-0.4500509424276378
9,13,-0.4583244495613092
9,14,-0.4865818826308876
9,15,-0.100009932057743
9,16,-0.08273351073040391
9,17,0.07468705617190707
9,18,0.07468705617190707
9,19,0.07468705617190707
9,20,0.07468705617190707
9,21,0.07468705617190707
9,22,0.07468705617190707
9,23,0.07468705617190707
9,24,0.07468705617190707
9,25,0.07468705617190707
9,26,0.07468705617190707
9,27,0.07468705617190707
9,28,0.07468705617190707
9,29,0.07468705617190707
9,30,0.07468705617190707
9,31,0.07468705617190707
9,32,0.07468705617190707
9,33,0.07468705617190707
9,34,0.07468705617190707
9,35,0.07468705617190707
9,36,0.07468705617190707
9,37,0.07468705617190707
9,38,0.07468705617190707
9,39,0.07468705617190707
9,40,0.07468705617190707
9,41,0.07468705617190707
9,42,0.07468705617190707
9,43,0.07468705617190707
9,44,0.07468705617190707
9,45,0.07468705617190707
9,46,0.07468705617190707
9,47,0.07468705617190707
9,48,0.07468705617190707
9,49,0.07468705617190707
9,50,0.07468705617190707
9,51,0.07468705617190707
9,52,0.07468705617190707
10,1,42735
10,14,-0.5189868785143603
10,15,-0.2397840863870959
10,16,-0.01252079199893785
10,17,0.3000376062062493
10,18,-0.1151820594382569
10,19,-0.05493014849097255
10,20,-0.05493014849097255
10,21,-0.05493014849097255
10,22,-0.05493014849097255
10,23,-0.05493014849097255
10,24,-0.05493014849097255
10,25,-0.05493014849097255
10,26,-0.05493014849097255
10,27,-0.05493014849097255
10,28,-0.05493014849097255
10,29,-0.05493014849097255
10,30,-0.05493014849097255
10,31,-0.05493014849097255
10,32,-0.05493014849097255
10,33,-0.05493014849097255
10,34,-0.05493014849097255
10,35,-0.05493014849097255
10,36,-0.05493014849097255
10,37,-0.05493014849097255
10,38,-0.05493014849097255
10,39,-0.05493014849097255
10,40,-0.05493014849097255
10,41,-0.05493014849097255
10,42,-0.05493014849097255
10,43,-0.05493014849097255
10,44,-0.05493014849097255
10,45,-0.05493014849097255
10,46,-0.05493014849097255
10,47,-0.05493014849097255
10,48,-0.05493014849097255
10,49,-0.05493014849097255
10,50,-0.05493014849097255
10,51,-0.05493014849097255
10,52,-0.05493014849097255
11,1,43100
11,16,-0.03999400040000101
11,17,0.2001500500062203
11,18,0.07348980370169844
11,19,0.1740313431290996
11,20,0.157394256377752
11,21,0.2820931576894115
11,22,0.2820931576894115
11,23,0.2820931576894115
11,24,0.2820931576894115
11,25,0.2820931576894115
11,26,0.2820931576894115
11,27,0.2820931576894115
11,28,0.2820931576894115
11,29,0.2820931576894115
11,30,0.2820931576894115
11,31,0.2820931576894115
11,32,0.2820931576894115
11,33,0.2820931576894115
11,34,0.2820931576894115
11,35,0.2820931576894115
11,36,0.2820931576894115
11,37,0.2820931576894115
11,38,0.2820931576894115
11,39,0.2820931576894115
11,40,0.2820931576894115
11,41,0.2820931576894115
11,42,0.2820931576894115
11,43,0.2820931576894115
11,44,0.2820931576894115
11,45,0.2820931576894115
11,46,0.2820931576894115
11,47,0.2820931576894115
11,48,0.2820931576894115
11,49,0.2820931576894115
11,50,0.2820931576894115
11,51,0.2820931576894115
11,52,0.2820931576894115
12,1,43465
12,18,0.07838989163155841
12,19,0.07838989163158061
12,20,0.06194937150048041
12,21,0.3390041783450259
12,22,0.2895071529679827
12,23,0.3426151435189873
12,24,0.187607693984293
12,25,0.2343541283920114
12,26,0.2343541283920114
12,27,0.2343541283920114
12,28,0.2343541283920114
12,29,0.2343541283920114
12,30,0.2343541283920114
12,31,0.2343541283920114
12,32,0.2343541283920114
12,33,0.2343541283920114
12,34,0.2343541283920114
12,35,0.2343541283920114
12,36,0.2343541283920114
12,37,0.2343541283920114
12,38,0.2343541283920114
12,39,0.2343541283920114
12,40,0.2343541283920114
12,41,0.2343541283920114
12,42,0.2343541283920114
12,43,0.2343541283920114
12,44,0.2343541283920114
12,45,0.2343541283920114
12,46,0.2343541283920114
12,47,0.2343541283920114
12,48,0.2343541283920114
12,49,0.2343541283920114
12,50,0.2343541283920114
12,51,0.2343541283920114
12,52,0.2343541283920114
13,1,43830
13,20,0.05248803135060598
13,21,0.2214516015410783
13,22,0.1190486358061627
13,23,0.1686730364466316
13,24,0.02358181985058216
13,25,0.2210188332817387
13,26,0.3712693419885671
13,27,0.1598952850611068
13,28,-0.02097628618118463
13,29,-0.009430310228020211
13,30,-0.009430310228020211
13,31,-0.009430310228020211
13,32,-0.009430310228020211
13,33,-0.009430310228020211
13,34,-0.009430310228020211
13,35,-0.009430310228020211
13,36,-0.009430310228020211
13,37,-0.009430310228020211
13,38,-0.009430310228020211
13,39,-0.009430310228020211
13,40,-0.009430310228020211
13,41,-0.009430310228020211
13,42,-0.009430310228020211
13,43,-0.009430310228020211
13,44,-0.009430310228020211
13,45,-0.009430310228020211
13,46,-0.009430310228020211
13,47,-0.009430310228020211
13,48,-0.009430310228020211
13,49,-0.009430310228020211
13,50,-0.009430310228020211
13,51,-0.009430310228020211
13,52,-0.009430310228020211
14,1,44196
14,23,0.1555800062495782
14,24,0.1685648701203801
14,25,0.2603046921112462
14,26,0.3003887663412641
14,27,0.132465972367557
14,28,-0.1339126064348495
14,29,-0.09571633453315798
14,30,-0.03054415496863694
14,31,-0.4923796969465988
14,32,-2.657403949513992
14,33,-2.657403949513992
14,34,-2.657403949513992
14,35,-2.657403949513992
14,36,-2.657403949513992
14,37,-2.657403949513992
14,38,-2.657403949513992
14,39,-2.657403949513992
14,40,-2.657403949513992
14,41,-2.657403949513992
14,42,-2.657403949513992
14,43,-2.657403949513992
14,44,-2.657403949513992
14,45,-2.657403949513992
14,46,-2.657403949513992
14,47,-2.657403949513992
14,48,-2.657403949513992
14,49,-2.657403949513992
14,50,-2.657403949513992
14,51,-2.657403949513992
14,52,-2.657403949513992
15,1,44561
15,27,0.1889333341656085
15,28,0.07509714884945673
15,29,0.08073518189748441
15,30,0.035598638033707
15,31,-0.3459257698102514
15,32,-2.423328265806446
15,33,-1.49562970548649
15,34,-0.2885033948250459
15,35,-0.5121403324772844
15,36,-0.3096364143617802
15,37,-0.3096364143617802
15,38,-0.3096364143617802
15,39,-0.3096364143617802
15,40,-0.3096364143617802
15,41,-0.3096364143617802
15,42,-0.3096364143617802
15,43,-0.3096364143617802
15,44,-0.3096364143617802
15,45,-0.3096364143617802
15,46,-0.3096364143617802
15,47,-0.3096364143617802
15,48,-0.3096364143617802
15,49,-0.3096364143617802
15,50,-0.3096364143617802
15,51,-0.3096364143617802
15,52,-0.3096364143617802
16,1,44926
16,31,-0.2239594630099373
16,32,-1.616318061533883
16,33,-1.337223818620836
16,34,-0.599348850912329
16,35,-1.073589070820447
16,36,-0.3018961902350958
16,37,-0.1048501255800471
16,38,-0.2588455356339781
16,39,-0.2454721753057276
16,40,-0.1730430455425092
16,41,-0.1730430455425092
16,42,-0.1730430455425092
16,43,-0.1730430455425092
16,44,-0.1730430455425092
16,45,-0.1730430455425092
16,46,-0.1730430455425092
16,47,-0.1730430455425092
16,48,-0.1730430455425092
16,49,-0.1730430455425092
16,50,-0.1730430455425092
16,51,-0.1730430455425092
16,52,-0.1730430455425092
17,1,45291
17,34,-0.9533114413926458
17,35,-1.288136903730974
17,36,-1.02314941214694
17,37,-0.9401582880721127
17,38,-1.151120647939763
17,39,-1.238905350026021
17,40,-0.8943276391025989
17,41,0.9692952624595019
17,42,0.8644693227634503
17,43,0.7038634017465073
17,44,0.6376744206510576
17,45,0.6376744206510576
17,46,0.6376744206510576
17,47,0.6376744206510576
17,48,0.6376744206510576
17,49,0.6376744206510576
17,50,0.6376744206510576
17,51,0.6376744206510576
17,52,0.6376744206510576
18,1,45657
18,38,-1.122870469184911
18,39,-1.257828108304415
18,40,-1.052658082693458
18,41,1.385373425334802
18,42,0.7652102000489602
18,43,0.4406734233171727
18,44,-0.09128981027868299
18,45,0.1544084105021826
18,46,0.001611361207976003
18,47,0.20168190406884
18,48,0.1856341247700399
18,49,0.1856341247700399
18,50,0.1856341247700399
18,51,0.1856341247700399
18,52,0.1856341247700399
19,1,46022
19,42,0.7168164388559273
19,43,0.4973302237080146
19,44,-0.0379413063031464
19,45,0.2685433396315773
19,46,0.01799217181808199
19,47,0.3613321345859122
19,48,0.4108497965175983
19,49,0.3997355152047577
19,50,0.08117592553187336
19,51,-0.06418790329880686
19,52,-0.09450306168263811
20,1,46387
20,46,0.04041576823738957
20,47,0.2315503018970322
20,48,0.3875789231538196
20,49,0.5005566802541939
20,50,0.3069836986764551
20,51,-0.08988642825158433
20,52,-0.3000102673190841
21,1,46752
21,50,0.4008103346141656
21,51,0.06556754296873635
21,52,-0.1252661784341358
22,1,47118
"@

$count = 0
foreach ($line in ($cellValues -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
    $count = $count + 1
}

Write-Output "Applied naive forecaster bugfix: $count cells written"
